# Memphis roster: three pairs of rows had their Player..bbref-url data
# (columns B..K) swapped while the No. column (A) stayed put:
#   row 6  <-> row 7   (Brandon Clarke          <-> David Roddy)
#   row 8  <-> row 9   (Ja Morant                <-> Jaren Jackson Jr.)
#   row 17 <-> row 18  (Vince Williams Jr. (TW)  <-> Luke Kennard)
#
# Use Range.Copy through a scratch row (20, just beyond the used range) so
# cell types/styles (e.g. the text-typed "Exp" column) are preserved exactly
# as Excel would when moving/sorting rows, instead of being re-inferred by
# assigning literal .Value strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $scratch = "B20:K20"
    $ws.Range("B${r1}:K${r1}").Copy($ws.Range($scratch))
    $ws.Range("B${r2}:K${r2}").Copy($ws.Range("B${r1}:K${r1}"))
    $ws.Range($scratch).Copy($ws.Range("B${r2}:K${r2}"))
    $ws.Range($scratch).Clear()
}

Swap-Rows 6 7
Swap-Rows 8 9
Swap-Rows 17 18
